# Auto-generated edit script: updates market-price derived values (H..N)
# on sheets ALC, ARM, BSM, CUL, GSM, LTW, WVR to match the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2878.52
$ws.Range("I19").Value = 3235.75
$ws.Range("K19").Value = 3235.75
$ws.Range("M19").Value = -3060.75

$ws.Range("H70").Value = 16515.584
$ws.Range("I70").Value = 1242.7142
$ws.Range("J70").Value = 37897.6
$ws.Range("K70").Value = 3728.1426
$ws.Range("L70").Value = 113692.8
$ws.Range("M70").Value = -3458.1426
$ws.Range("N70").Value = -114232.8

$ws.Range("H73").Value = 16515.584
$ws.Range("I73").Value = 1242.7142
$ws.Range("J73").Value = 37897.6
$ws.Range("K73").Value = 3728.1426
$ws.Range("L73").Value = 113692.8
$ws.Range("M73").Value = -2792.1426
$ws.Range("N73").Value = -115564.8

$ws.Range("H111").Value = 4404.5
$ws.Range("I111").Value = 3329.1428
$ws.Range("J111").Value = 11932
$ws.Range("K111").Value = 9987.428400000001
$ws.Range("L111").Value = 35796
$ws.Range("M111").Value = -6920.428400000001
$ws.Range("N111").Value = -41930

$ws.Range("H132").Value = 38646.723
$ws.Range("I132").Value = 44913.668
$ws.Range("J132").Value = 7881.727
$ws.Range("K132").Value = 134741.004
$ws.Range("L132").Value = 23645.181
$ws.Range("M132").Value = -132211.004
$ws.Range("N132").Value = -28705.181

$ws.Range("H137").Value = 38240876
$ws.Range("I137").Value = 83334490
$ws.Range("J137").Value = 2165984.5
$ws.Range("K137").Value = 250003470
$ws.Range("L137").Value = 6497953.5
$ws.Range("M137").Value = -250000920
$ws.Range("N137").Value = -6503053.5

$ws.Range("H138").Value = 2273.4902
$ws.Range("I138").Value = 1301.8125
$ws.Range("J138").Value = 2717.6858
$ws.Range("K138").Value = 3905.4375
$ws.Range("L138").Value = 8153.057400000001
$ws.Range("M138").Value = 1234.5625
$ws.Range("N138").Value = -18433.0574

$ws.Range("H141").Value = 2114.3333
$ws.Range("I141").Value = 2719
$ws.Range("K141").Value = 8157
$ws.Range("M141").Value = -2977

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6250824.5
$ws.Range("I32").Value = 6329937.5
$ws.Range("K32").Value = 6329937.5
$ws.Range("M32").Value = -6329650.5

$ws.Range("H61").Value = 1044283.2
$ws.Range("I61").Value = 1588916.2
$ws.Range("J61").Value = 4529.273
$ws.Range("K61").Value = 1588916.2
$ws.Range("L61").Value = 4529.273
$ws.Range("M61").Value = -1588704.2
$ws.Range("N61").Value = -4953.273

$ws.Range("H110").Value = 1727.25
$ws.Range("I110").Value = 1897.5714
$ws.Range("J110").Value = 1488.8
$ws.Range("K110").Value = 1897.5714
$ws.Range("L110").Value = 1488.8
$ws.Range("M110").Value = 147.4286
$ws.Range("N110").Value = -5578.8

$ws.Range("H122").Value = 1703.5
$ws.Range("I122").Value = 1161.2142
$ws.Range("K122").Value = 3483.6426
$ws.Range("M122").Value = -1033.6426

$ws.Range("H136").Value = 1044283.2
$ws.Range("I136").Value = 1588916.2
$ws.Range("J136").Value = 4529.273
$ws.Range("K136").Value = 4766748.6
$ws.Range("L136").Value = 13587.819
$ws.Range("M136").Value = -4764198.6
$ws.Range("N136").Value = -18687.819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4553.2
$ws.Range("I107").Value = 4522.923
$ws.Range("K107").Value = 4522.923
$ws.Range("M107").Value = -2602.923

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()

$ws.Range("H117").Value = 19611784
$ws.Range("J117").Value = 4829.154
$ws.Range("L117").Value = 14487.462
$ws.Range("N117").Value = -21371.462

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("N118").ClearContents()

$ws.Range("H119").Value = 704.5
$ws.Range("I119").Value = 704.5
$ws.Range("K119").Value = 2113.5
$ws.Range("M119").Value = 2724.5

$ws.Range("H120").Value = 23332.75
$ws.Range("I120").Value = 11632.5
$ws.Range("K120").Value = 34897.5
$ws.Range("M120").Value = -30059.5

$ws.Range("H121").Value = 66668080
$ws.Range("I121").Value = 100000120
$ws.Range("J121").Value = 4000
$ws.Range("K121").Value = 300000360
$ws.Range("L121").Value = 12000
$ws.Range("M121").Value = -299999050
$ws.Range("N121").Value = -14620

$ws.Range("H125").Value = 18693.666
$ws.Range("J125").Value = 25033
$ws.Range("L125").Value = 75099
$ws.Range("N125").Value = -84939

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 10314
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 507.55
$ws.Range("I16").Value = 504
$ws.Range("J16").Value = 511.8889
$ws.Range("K16").Value = 504
$ws.Range("L16").Value = 511.8889
$ws.Range("M16").Value = -334
$ws.Range("N16").Value = -851.8888999999999

$ws.Range("H61").Value = 2400.2
$ws.Range("I61").Value = 1501.5
$ws.Range("K61").Value = 1501.5
$ws.Range("M61").Value = -1299.5

$ws.Range("H93").Value = 1751
$ws.Range("I93").Value = 1790.2222
$ws.Range("J93").Value = 1574.5
$ws.Range("K93").Value = 1790.2222
$ws.Range("L93").Value = 1574.5
$ws.Range("M93").Value = -542.2221999999999
$ws.Range("N93").Value = -4070.5

$ws.Range("H113").Value = 2400.2
$ws.Range("I113").Value = 1501.5
$ws.Range("K113").Value = 1501.5
$ws.Range("M113").Value = 668.5

$ws.Range("H122").Value = 5293.8237
$ws.Range("I122").Value = 4768.846
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 14306.538
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -11856.538
$ws.Range("N122").Value = -25900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3411.5806
$ws.Range("I122").Value = 3275.4666
$ws.Range("K122").Value = 9826.399800000001
$ws.Range("M122").Value = -7376.399800000001

$ws.Range("H126").Value = 4970.143
$ws.Range("I126").Value = 4134.3335
$ws.Range("K126").Value = 12403.0005
$ws.Range("M126").Value = -9933.000499999998
